$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# "Sprint No." value cell: change "1" -> "2"
$sprintCell = $t.Cell(2, 4)
$sprintRange = $sprintCell.Range
$sprintRange.MoveEnd(1, -1)
$sprintRange.Text = "2"

# "Review Date" value cell: change "02/09/18" -> "02/21/18"
$dateCell = $t.Cell(3, 2)
$dateRange = $dateCell.Range
$dateRange.MoveEnd(1, -1)
$dateRange.Text = "02/21/18"
